$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the previously used range (A1:F6) completely so stale
# formatting/values in column F and the old A/B/C layout are gone.
$ws.Range("A1:F6").Clear()

# --- Header row (row 1): styled like the old header cells (bold font,
# thin border, centered horizontal / top vertical alignment) ---
$headers = @("EL_Astral25", "FNRATE_EXACT_ASTRAL", "TAXON", "MODELCONDITION", "GENE")
for ($col = 1; $col -le 5; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# --- Data rows 2-6: A=number, B=0, C="11-texon", D="simulated_25genes_strongILS", E=number ---
$dataA = @(114, 29, 29, 114, 29)
$dataE = @(3, 8, 14, 16, 18)

for ($i = 0; $i -lt 5; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $dataA[$i]
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = "11-texon"
    $ws.Cells.Item($r, 4).Value = "simulated_25genes_strongILS"
    $ws.Cells.Item($r, 5).Value = $dataE[$i]
}
